$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently sits right
#    after the title (Heading1) paragraph.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Insert a new paragraph before the closing "Prompt: ..." paragraph,
#    containing the bold title text "Play Solar Disc Slot Game for Free -
#    Review and Recommendations".
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($count)
$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Solar Disc Slot Game for Free - Review and Recommendations</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($titleXml)

# 3. Replace the text of the old "Prompt: ..." paragraph (now the last
#    paragraph again) with the meta-description text, keeping its italic
#    formatting intact.
$oldPrompt = "Prompt: Create a cartoon-style feature image for the game " + [char]34 + "Solar Disc" + [char]34 + " featuring a happy Maya warrior wearing glasses. The image should be eye-catching and highlight the Aztec theme of the game, with symbols like coins, keys, chalices, and the powerful sun god Huitzilopochtli. The Maya warrior should be front and center, with a big smile on their face, wearing traditional Aztec armor and wielding a spear or other weapon. The background can be a desert landscape with the Solar Disc symbol shining bright in the sky. Make sure the image includes the game title and any relevant branding or logos. The overall feel should be fun and adventurous, inspiring players to join the Maya warrior on an exciting journey to win big with Solar Disc."
$newDescription = "Read a review of Solar Disc slot game, with pros and cons. Play for free and discover similar Aztec-themed slot games such as Aztec Gems."

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2)
